# The checklist's "icon" column (A) used to hold a single-space placeholder
# string in each row; now that a real icon renders in its place, the
# placeholder text is no longer needed, so blank out A1:A5 while leaving the
# True/False flags in column B untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 5; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.ClearContents()
    # Touch a format property (already at its default) so the cell keeps a
    # presence on the sheet (and the used range stays A1:B5) without
    # actually altering its appearance/style.
    $cell.Font.Bold = $false
}
